# KIBON-120: Benutzer statistik
# Replace hard-coded German header labels with "{...Title}" placeholders that
# get resolved/translated at report-generation time (so the report can be
# produced in the user's language instead of only German).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column header row (row 5) - swap the static German captions for title
# placeholders so they can be localized. (Filled in the same left-to-right
# order as the original edit, with the N/O pair set in reverse.)
$ws.Range("A5").Value = "{usernameTitle}"
$ws.Range("B5").Value = "{vornameTitle}"
$ws.Range("C5").Value = "{nachnameTitle}"
$ws.Range("D5").Value = "{emailTitle}"
$ws.Range("E5").Value = "{roleTitle}"
$ws.Range("F5").Value = "{roleGueltigBisTitel}"
$ws.Range("G5").Value = "{gemeindenTitle}"
$ws.Range("H5").Value = "{institutionTitle}"
$ws.Range("I5").Value = "{traegerschaftTitle}"
$ws.Range("J5").Value = "{kitaTitel}"
$ws.Range("K5").Value = "{tagesfamilienTitle}"
$ws.Range("L5").Value = "{tagesschulenTitel}"
$ws.Range("M5").Value = "{ferieninselTitle}"
$ws.Range("O5").Value = "{isSchulamtTitle}"
$ws.Range("N5").Value = "{isJugendamtTitle}"
$ws.Range("P5").Value = "{statusTitle}"

# Row 6 (the {placeholder} row used by the report engine to repeat data rows)
# is unchanged in content - left as-is.

# "Stichtag" label (row 3) - B3 already holds the {stichtag} placeholder and
# stays untouched.
$ws.Range("A3").Value = "{stichtagTitle}"

# Report title (row 1)
$ws.Range("A1").Value = "{reportBenutzerTitle}"

# Column width tweaks that came along with the edit in Excel.
$ws.Columns.Item(6).ColumnWidth = 16
$ws.Columns.Item(10).ColumnWidth = 12.6666666666667

# The saved sheet view used to pin an odd full-column selection
# (L1:L1048576); reset it back to the default single-cell selection.
$ws.Range("A1").Select()
